# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn and de-de sheets to
# reflect that the file is now "Ready for handoff" (instead of the old
# "Handed back: in sync with en-US" status), with refreshed handoff
# metadata (new handoff file names / datetimes) and an error detail
# message about the handback file being out of date.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/707257f2ab2b7a3a718550fbc904aa009b823fd6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/879a334c6f6dd5d4ca7d1cb9c265716dca7d7180/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-24 09:19:00"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").ClearFormats()
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-24 09:18:48"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").ClearFormats()
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-24 09:19:00"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
